$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C32").Value = 50.55
$ws.Range("D33").Value = 66.5
$ws.Range("C59").Value = 13.35
$ws.Range("C94").Value = 13.35
$ws.Range("C114").Value = 13.35
$ws.Range("C128").Value = 13.35
$ws.Range("C130").Value = 50.55
